$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 3668
$ws.Range("I7").Value = 4005
$ws.Range("J7").Value = 3499.5
$ws.Range("K7").Value = 4005
$ws.Range("L7").Value = 3499.5
$ws.Range("M7").Value = -3893
$ws.Range("N7").Value = -3723.5
$ws.Range("H14").Value = 3668
$ws.Range("I14").Value = 4005
$ws.Range("J14").Value = 3499.5
$ws.Range("K14").Value = 4005
$ws.Range("L14").Value = 3499.5
$ws.Range("M14").Value = -3814
$ws.Range("N14").Value = -3881.5
$ws.Range("H40").Value = 2574.5
$ws.Range("J40").Value = 2599.3333
$ws.Range("L40").Value = 2599.3333
$ws.Range("N40").Value = -2949.3333
$ws.Range("H92").Value = 1803
$ws.Range("J92").Value = 4999
$ws.Range("L92").Value = 4999
$ws.Range("N92").Value = -7495
$ws.Range("H96").Value = 4000
$ws.Range("J96").Value = 4000
$ws.Range("L96").Value = 12000
$ws.Range("N96").Value = -14746
$ws.Range("H98").Value = 931.0909
$ws.Range("I98").Value = 924.2
$ws.Range("K98").Value = 924.2
$ws.Range("M98").Value = 573.8
$ws.Range("H101").Value = 2480.8572
$ws.Range("I101").Value = 2773.2
$ws.Range("J101").Value = 1750
$ws.Range("K101").Value = 8319.599999999999
$ws.Range("L101").Value = 5250
$ws.Range("M101").Value = -6697.599999999999
$ws.Range("N101").Value = -8494
$ws.Range("H122").Value = 931.0909
$ws.Range("I122").Value = 924.2
$ws.Range("K122").Value = 2772.6
$ws.Range("M122").Value = -322.6000000000004
$ws.Range("H137").Value = 2542.7144
$ws.Range("I137").Value = 2500
$ws.Range("K137").Value = 7500
$ws.Range("M137").Value = -4950
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 11289.2
$ws.Range("J29").Value = 11289.2
$ws.Range("L29").Value = 11289.2
$ws.Range("N29").Value = -11905.2
$ws.Range("H41").Value = 5401.25
$ws.Range("I41").Value = 1302.5
$ws.Range("K41").Value = 1302.5
$ws.Range("M41").Value = -888.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H134").Value = 1762.4
$ws.Range("I134").Value = 1624.8889
$ws.Range("K134").Value = 4874.6667
$ws.Range("M134").Value = -2339.6667
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 151.4
$ws.Range("I6").Value = 86
$ws.Range("J6").Value = 195
$ws.Range("K6").Value = 86
$ws.Range("L6").Value = 195
$ws.Range("M6").Value = 27
$ws.Range("N6").Value = -421
$ws.Range("H23").Value = 10666.667
$ws.Range("J23").Value = 10666.667
$ws.Range("L23").Value = 10666.667
$ws.Range("N23").Value = -11146.667
$ws.Range("H27").Value = 10666.667
$ws.Range("J27").Value = 10666.667
$ws.Range("L27").Value = 10666.667
$ws.Range("N27").Value = -11050.667
$ws.Range("H31").Value = 1861.8182
$ws.Range("I31").Value = 1663
$ws.Range("J31").Value = 2756.5
$ws.Range("K31").Value = 1663
$ws.Range("L31").Value = 2756.5
$ws.Range("M31").Value = -1368
$ws.Range("N31").Value = -3346.5
$ws.Range("H34").Value = 1861.8182
$ws.Range("I34").Value = 1663
$ws.Range("J34").Value = 2756.5
$ws.Range("K34").Value = 1663
$ws.Range("L34").Value = 2756.5
$ws.Range("M34").Value = -1461
$ws.Range("N34").Value = -3160.5
$ws.Range("H58").Value = 2042.4
$ws.Range("I58").Value = 1303
$ws.Range("K58").Value = 1303
$ws.Range("M58").Value = -1100
$ws.Range("H106").Value = 28125
$ws.Range("J106").Value = 28125
$ws.Range("L106").Value = 28125
$ws.Range("N106").Value = -30649
$ws.Range("H136").Value = 2042.4
$ws.Range("I136").Value = 1303
$ws.Range("K136").Value = 3909
$ws.Range("M136").Value = -1359
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 117647610
$ws.Range("I4").Value = 181818660
$ws.Range("J4").Value = 686.5
$ws.Range("K4").Value = 545455980
$ws.Range("L4").Value = 2059.5
$ws.Range("M4").Value = -545455868
$ws.Range("N4").Value = -2283.5
$ws.Range("H7").Value = 309
$ws.Range("I7").Value = 224.66667
$ws.Range("J7").Value = 345.14285
$ws.Range("K7").Value = 674.00001
$ws.Range("L7").Value = 1035.42855
$ws.Range("M7").Value = -562.00001
$ws.Range("N7").Value = -1259.42855
$ws.Range("H13").Value = 60.375
$ws.Range("I13").Value = 25.5
$ws.Range("J13").Value = 95.25
$ws.Range("K13").Value = 76.5
$ws.Range("L13").Value = 285.75
$ws.Range("M13").Value = 91.5
$ws.Range("N13").Value = -621.75
$ws.Range("H108").Value = 688.3333
$ws.Range("I108").Value = 688.3333
$ws.Range("K108").Value = 2064.9999
$ws.Range("M108").Value = 815.0001000000002
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H117").Value = 677
$ws.Range("I117").Value = 499.5
$ws.Range("J117").Value = 1032
$ws.Range("K117").Value = 1498.5
$ws.Range("L117").Value = 3096
$ws.Range("M117").Value = 1943.5
$ws.Range("N117").Value = -9980
$ws.Range("H132").Value = 1259.3334
$ws.Range("I132").Value = 999
$ws.Range("K132").Value = 8991
$ws.Range("M132").Value = -6461
$ws.Range("H137").Value = 1326.6666
$ws.Range("I137").Value = 1240
$ws.Range("J137").Value = 1500
$ws.Range("K137").Value = 3720
$ws.Range("L137").Value = 4500
$ws.Range("M137").Value = 1380
$ws.Range("N137").Value = -14700
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1500.6666
$ws.Range("I102").Value = 1251.5
$ws.Range("J102").Value = 1999
$ws.Range("K102").Value = 1251.5
$ws.Range("L102").Value = 1999
$ws.Range("M102").Value = 370.5
$ws.Range("N102").Value = -5243
$ws.Range("H113").Value = 1719.7
$ws.Range("J113").Value = 1972.5
$ws.Range("L113").Value = 1972.5
$ws.Range("N113").Value = -6312.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31352
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3548.875
$ws.Range("I126").Value = 3548.875
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10646.625
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -8176.625
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 3476.25
$ws.Range("I132").Value = 2900
$ws.Range("J132").Value = 3668.3333
$ws.Range("K132").Value = 8700
$ws.Range("L132").Value = 11004.9999
$ws.Range("M132").Value = -6170
$ws.Range("N132").Value = -16064.9999
